# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6811
$wsExhibit.Range("F12").Value = 21
$wsExhibit.Range("F13").Value = 177
$wsExhibit.Range("F18").Value = 3431
$wsExhibit.Range("F19").Value = 20
$wsExhibit.Range("F22").Value = 2064
$wsExhibit.Range("F23").Value = 163
$wsExhibit.Range("F24").Value = 3
$wsExhibit.Range("F28").Value = 10

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6811
$wsAll.Range("F13").Value = 21
$wsAll.Range("F14").Value = 177
$wsAll.Range("F19").Value = 3431
$wsAll.Range("F20").Value = 20
$wsAll.Range("F23").Value = 2064
$wsAll.Range("F24").Value = 163
$wsAll.Range("F25").Value = 3
$wsAll.Range("F29").Value = 10
